$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value
$updates = @(
    @('D2', '44.689.18'),
    @('E2', '  +4.07%  '),
    @('D3', '2.423.69'),
    @('E3', '  +2.57%  '),
    @('D4', '1.00'),
    @('E4', '  -0.02%  '),
    @('D5', '315.91'),
    @('E5', '  +4.15%  '),
    @('D6', '101.95'),
    @('E6', '  +7.24%  '),
    @('E7', '  +2.53%  '),
    @('E8', '  -0.06%  '),
    @('D9', '0.525'),
    @('E9', '  +9.09%  '),
    @('D10', '35.57'),
    @('E10', '  +4.39%  '),
    @('D11', '0.0800'),
    @('E11', '  +2.00%  '),
    @('D12', '18.93'),
    @('E12', '  +2.61%  '),
    @('E13', '  -2.20%  '),
    @('D14', '6.96'),
    @('E14', '  +3.52%  '),
    @('D15', '2.802.86'),
    @('E15', '  +2.72%  '),
    @('D16', '2.374.15'),
    @('E16', '  +1.87%  '),
    @('D17', '0.833'),
    @('E17', '  +4.95%  '),
    @('D18', '44.560.12'),
    @('E18', '  +3.84%  '),
    @('D19', '12.37'),
    @('E19', '  +3.85%  '),
    @('D20', '6.40'),
    @('E20', '  +2.30%  '),
    @('D21', '0.0₃0921'),
    @('E21', '  +4.17%  '),
    @('D22', '68.83'),
    @('E22', '  +1.08%  '),
    @('D23', '242.51'),
    @('E23', '  +3.22%  '),
    @('E24', '  +5.89%  '),
    @('E25', '  +1.51%  '),
    @('E26', '  +0.05%  '),
    @('D27', '25.25'),
    @('E27', '  +3.10%  '),
    @('D28', '2.27'),
    @('E28', '  -4.03%  '),
    @('D29', '9.51'),
    @('E29', '  +1.78%  '),
    @('D30', '33.61'),
    @('E30', '  +5.05%  '),
    @('D31', '48.42'),
    @('E31', '  +1.28%  '),
    @('E32', '  +20.50%  '),
    @('D33', '19.43'),
    @('E33', '  +10.79%  '),
    @('B34', 'Filecoin'),
    @('C34', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'),
    @('D34', '5.18'),
    @('E34', '  +3.76%  '),
    @('B35', 'Hedera'),
    @('C35', 'https://coinranking.com/coin/jad286TjB+hedera-hbar'),
    @('D35', '0.0776'),
    @('E35', '  +8.56%  '),
    @('E36', '  +0.26%  '),
    @('E37', '  +3.07%  '),
    @('E38', '  +3.99%  '),
    @('E39', '  +1.33%  '),
    @('B40', 'WEMIXToken'),
    @('C40', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'),
    @('D40', '2.21'),
    @('E40', '  -2.34%  '),
    @('B41', 'Monero'),
    @('C41', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'),
    @('D41', '120.07'),
    @('E41', '  -3.61%  '),
    @('B42', 'Stellar'),
    @('C42', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'),
    @('D42', '0.109'),
    @('E42', '  +1.95%  '),
    @('D43', '21.14'),
    @('E43', '  -1.48%  '),
    @('D44', '0.0290'),
    @('E44', '  +4.40%  '),
    @('D45', '1.943.64'),
    @('E45', '  +0.66%  '),
    @('D46', '2.16'),
    @('E46', '  +1.06%  '),
    @('E47', '  +9.20%  '),
    @('D48', '9.44'),
    @('E48', '  +1.93%  '),
    @('D49', '1.67'),
    @('E49', '  +10.95%  '),
    @('B50', 'BitcoinSV'),
    @('C50', 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'),
    @('D50', '75.54'),
    @('E50', '  +5.81%  '),
    @('B51', 'MultiversX'),
    @('C51', 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'),
    @('D51', '54.53'),
    @('E51', '  +6.53%  ')
)

foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    $cell = $ws.Range($ref)
    if ($ref.StartsWith("D")) {
        # Price column: force text so numeric-looking values (e.g. "1.00", "0.0800")
        # keep their exact formatting instead of being parsed as numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}
